# TC01_Trials_Filter_AssocFileFormat-Bai.xlsx - add StatQuery column with the
# bai-file-format trial/arm/case query, and the file-count summary query,
# inserted as new column B (existing dbExcel / WebExcel data shifts right).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$caseQuery = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s WHERE f.file_format IN ['bai']  RETURN DISTINCT coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"
$statQuery = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s , f WHERE f.file_format IN ['bai','bam','vcf'] RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Insert a new column before column B so the existing dbExcel/Neo4j/Web data
# shifts from B/C to C/D, making room for the new StatQuery column.
$ws.Columns("B").Insert()

$ws.Range("B1").Value = "StatQuery"
$ws.Range("A2").Value = $caseQuery
$ws.Range("B2").Value = $statQuery

$ws.Range("A2").WrapText = $true
$ws.Range("B2").WrapText = $true

$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth
$ws.Rows(2).RowHeight = 101.5

$ws.Columns("B").Select()
